$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First "Shop site" (the list item right after "Home") -> "Collection".
#    Word leaves its auto-tracked "_GoBack" bookmark (the "last edit" marker)
#    collapsed immediately after the edited text, so we reproduce that too.
#
#    We replace with a one-character-longer placeholder "CollectionX" so the
#    trailing "X" gives us a real (non-collapsed) run to anchor the bookmark
#    on; bookmarking a truly zero-width point turns out to land in the wrong
#    spot. We then add the bookmark around that "X" and delete it, which
#    leaves the bookmark collapsed right after "Collection", exactly where
#    Word would have left it.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Shop site", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "CollectionX", 1)

if ($r.Find.Found -or $r.Text -eq "CollectionX") {
    $markerStart = $r.End - 1
    $markerEnd = $r.End
    $bmRange = $d.Range($markerStart, $markerEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $delRange = $d.Range($markerStart, $markerEnd)
    $delRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Second "Shop site" (the standalone paragraph further down) -> "Collection".
#    No bookmark involved here.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Shop site", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "Collection", 1)

# ---------------------------------------------------------------------------
# Note: adding the new "_GoBack" bookmark above automatically removes the
# old "_GoBack" bookmark that used to sit after the final "Footer" paragraph
# (bookmark names are unique per document), matching the target diff.
# ---------------------------------------------------------------------------
